$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 15:57"

# Row 4
$ws.Range("B4").Value = 8038543
$ws.Range("C4").Value = 754
$ws.Range("D4").Value = 5186407
$ws.Range("E4").Value = 2632109
$ws.Range("G4").Value = 16
$ws.Range("H4").Value = 220027

# Row 5
$ws.Range("B5").Value = 7179006
$ws.Range("C5").Value = 5441
$ws.Range("D5").Value = 6228098
$ws.Range("E5").Value = 840985
$ws.Range("G5").Value = 29
$ws.Range("H5").Value = 109923

# Row 18
$ws.Range("B18").Value = 409358
$ws.Range("C18").Value = 3921
$ws.Range("D18").Value = 344208
$ws.Range("E18").Value = 55180
$ws.Range("G18").Value = 58
$ws.Range("H18").Value = 9970

# Row 23
$ws.Range("B23").Value = 340089
$ws.Range("C23").Value = 474
$ws.Range("D23").Value = 326339
$ws.Range("E23").Value = 8663
$ws.Range("G23").Value = 19
$ws.Range("H23").Value = 5087

# Row 25
$ws.Range("B25").Value = 332444
$ws.Range("C25").Value = 1350
$ws.Range("E25").Value = 43617
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 9727

# Row 51
$ws.Range("B51").Value = 89121
$ws.Range("C51").Value = 1208
$ws.Range("D51").Value = 54047
$ws.Range("E51").Value = 32964
$ws.Range("G51").Value = 16
$ws.Range("H51").Value = 2110

# Row 57
$ws.Range("E57").Value = 3977
$ws.Range("G57").Value = 4
$ws.Range("H57").Value = 284

# Row 58
$ws.Range("D58").Value = 49800
$ws.Range("E58").Value = 13979
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 2102

# Row 60
$ws.Range("B60").Value = 61642
$ws.Range("C60").Value = 323
$ws.Range("D60").Value = 58613
$ws.Range("E60").Value = 2518
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 511

# Row 69
$ws.Range("B69").Value = 47097
$ws.Range("C69").Value = 67
$ws.Range("D69").Value = 46452
$ws.Range("E69").Value = 337

# Row 74
$ws.Range("B74").Value = 41937
$ws.Range("C74").Value = 318
$ws.Range("D74").Value = 31340
$ws.Range("E74").Value = 9810
$ws.Range("G74").Value = 10
$ws.Range("H74").Value = 787

# Row 77
$ws.Range("B77").Value = 35006
$ws.Range("C77").Value = 152
$ws.Range("E77").Value = 2703
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = 767

# Row 80
$ws.Range("B80").Value = 31173
$ws.Range("C80").Value = 336
$ws.Range("D80").Value = 23777
$ws.Range("E80").Value = 6449
$ws.Range("G80").Value = 11
$ws.Range("H80").Value = 947

# Row 81
$ws.Range("C81").Value = 145
$ws.Range("G81").Value = 5

# Row 82
$ws.Range("B82").Value = 30437
$ws.Range("C82").Value = 1123
$ws.Range("D82").Value = 11769
$ws.Range("E82").Value = 17975
$ws.Range("G82").Value = 29
$ws.Range("H82").Value = 693

# Row 95
$ws.Range("B95").Value = 15730
$ws.Range("C95").Value = 91
$ws.Range("E95").Value = 3590

# Row 106
$ws.Range("B106").Value = 10872
$ws.Range("C106").Value = 4
$ws.Range("E106").Value = 347

# Row 107
$ws.Range("B107").Value = 10297
$ws.Range("C107").Value = 37
$ws.Range("D107").Value = 9177
$ws.Range("E107").Value = 1041
